# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record for Acelga (Femacal de La Calera) is inserted
# as row 336, pushing the existing rows 336:383 down to 337:384.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 336 (shifts rows 336:383 -> 337:384,
# also carries the dimension ref to A1:R384 and preserves the column D
# date-format style on the shifted cells).
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with the new record.
$ws.Cells.Item(336, 1).Value  = 3
$ws.Cells.Item(336, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(336, 3).Value  = "Coquimbo"
$ws.Cells.Item(336, 4).Value  = 44776
$ws.Cells.Item(336, 5).Value  = 5
$ws.Cells.Item(336, 6).Value  = 100112009
$ws.Cells.Item(336, 7).Value  = "Acelga"
$ws.Cells.Item(336, 8).Value  = "Sin especificar"
$ws.Cells.Item(336, 9).Value  = "Primera"
$ws.Cells.Item(336, 10).Value = 240
$ws.Cells.Item(336, 11).Value = 3300
$ws.Cells.Item(336, 12).Value = 3500
$ws.Cells.Item(336, 13).Value = 3400
$ws.Cells.Item(336, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(336, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(336, 16).Value = 567
$ws.Cells.Item(336, 17).Value = 6
$ws.Cells.Item(336, 18).Value = "Hortaliza"
